$d = $word.ActiveDocument

function Find-ParaContaining($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$text*") {
            return $p
        }
    }
    return $null
}

# --- Step 1: drop the "Menu déroulant qui prend le verrou." bullet
#     entirely; the following bullet ("Ouverture de tuile ... dirty.")
#     slides up into its numbered-list slot.
$pMenu = Find-ParaContaining("Menu d*roulant qui prend le verrou")
if ($pMenu -ne $null) {
    $pMenu.Range.Delete()
}

# --- Step 2: drop the "Un Resolver absent n'est pas signalé ... CoreData."
#     bullet entirely.
$pResolver = Find-ParaContaining("Un Resolver absent")
if ($pResolver -ne $null) {
    $pResolver.Range.Delete()
}

# --- Step 3: relocate the hidden "_GoBack" bookmark from the trailing
#     empty paragraph to right before the paragraph mark of the
#     "ConnectionString câblé en dur." bullet, leaving a genuinely empty
#     paragraph at the very end of the document.
$pConn = Find-ParaContaining("c*bl* en dur")
if ($pConn -ne $null) {
    $target = $pConn.Range.End - 1

    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()

    $insertRng = $d.Range($target, $target)
    $insertRng.InsertAfter("X")

    $xRng = $d.Range($target, $target + 1)
    $d.Bookmarks.Add("_GoBack", $xRng)

    $xRng2 = $d.Range($target, $target + 1)
    $xRng2.Text = ""
}
